$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing A40 timestamp (tiny floating point precision update)
$ws.Range("A40").Value = 44353.79480840857

# Append the new data row (row 41) retrieved in this run
$ws.Range("A41").Value = 44354.83177632051
$ws.Range("B41").Value = 73336
$ws.Range("C41").Value = 61868
$ws.Range("D41").Value = 3279
$ws.Range("E41").Value = 2094
$ws.Range("F41").Value = 1472
$ws.Range("G41").Value = 19256
$ws.Range("H41").Value = 1402
$ws.Range("I41").Value = 880
$ws.Range("J41").Value = 202

# Match the date-time style used by the rest of column A
$ws.Range("A41").NumberFormat = $ws.Range("A40").NumberFormat
